# Handback status report refresh: the in-flight file
# (8870f70d-9f1b-4ed2-befb-a238c465f730.md) was handed back and renamed on
# disk to 573f8fe1-0671-4552-bfd8-037bdf8d7374.md, and a brand-new file
# (ec2c564e-5368-466f-bc81-74a9e2afe9c5.md) entered the pipeline alongside
# it. Refresh all three report sheets (Overview, zh-cn, de-de) to reflect
# both: update the renamed-file row and append a new row for the new file.

$wb = $excel.ActiveWorkbook

$oldGuid = "8870f70d-9f1b-4ed2-befb-a238c465f730"
$renamedGuid = "573f8fe1-0671-4552-bfd8-037bdf8d7374"
$newGuid = "ec2c564e-5368-466f-bc81-74a9e2afe9c5"

$renamedHash = "e7a6bd797cab50c2f65adc9827b7c04bdeeccde0"
$newHash = "92cde374c60cfa040c8046cb03d2a078fd95248b"

$renamedFile = "$renamedGuid.md"
$newFile = "$newGuid.md"

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

# Drop every existing hyperlink on the sheet up front -- deleting one via
# Range.Hyperlinks.Delete() clears the whole sheet's collection in this
# host, so the only reliable sequencing is delete-everything-then-re-add.
$wsOverview.Range("A1").Hyperlinks.Delete()

# Row 2: rename in place, refresh the "Latest HO Xliff Generate Date".
$wsOverview.Range("A2").Value = $renamedFile
$wsOverview.Range("B2").Value = "e2e\$renamedFile"
$wsOverview.Range("C2").Value = ".md"
$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("G2").Value = "2016-08-29 21:15:56"

# Row 3: brand-new file.
$wsOverview.Range("A3").Value = $newFile
$wsOverview.Range("B3").Value = "e2e\$newFile"
$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("G3").Value = "2016-08-29 21:15:56"

$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/50f8611c0f7b5f3ff72b035e33299bbcfd184205/e2e/$renamedFile", "", "", "e2e\$renamedFile") | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/50f8611c0f7b5f3ff72b035e33299bbcfd184205/e2e/$newFile", "", "", "e2e\$newFile") | Out-Null

$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.Resize($wsOverview.Range("A1:G3"))

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("A1").Hyperlinks.Delete()

$zhRenamedXlf = "$renamedGuid.$renamedHash.zh-cn.xlf"
$zhNewXlf = "$newGuid.$newHash.zh-cn.xlf"

# Row 2: rename in place, refresh handoff/handback timestamps.
$wsZh.Range("A2").Value = $renamedFile
$wsZh.Range("B2").Value = ".md"
$wsZh.Range("C2").Value = "Handed back: in sync with en-US"
$wsZh.Range("D2").Value = "e2e"
$wsZh.Range("E2").Value = "ht"
$wsZh.Range("F2").Value = "False"
$wsZh.Range("G2").Value = $zhRenamedXlf
$wsZh.Range("H2").Value = "2016-08-29 21:15:51"
$wsZh.Range("I2").Value = $renamedFile
$wsZh.Range("J2").Value = $zhRenamedXlf
$wsZh.Range("K2").Value = "2016-08-29 21:16:19"
$wsZh.Range("L2").Value = ""
$wsZh.Range("M2").Value = "True"
$wsZh.Range("N2").Value = ""
$wsZh.Range("O2").Value = "False"
$wsZh.Range("P2").Value = ""

# Row 3: brand-new file.
$wsZh.Range("A3").Value = $newFile
$wsZh.Range("B3").Value = ".md"
$wsZh.Range("C3").Value = "Handed back: in sync with en-US"
$wsZh.Range("D3").Value = "e2e"
$wsZh.Range("E3").Value = "ht"
$wsZh.Range("F3").Value = "True"
$wsZh.Range("G3").Value = $zhNewXlf
$wsZh.Range("H3").Value = "2016-08-29 21:15:51"
$wsZh.Range("I3").Value = $newFile
$wsZh.Range("J3").Value = $zhNewXlf
$wsZh.Range("K3").Value = "2016-08-29 21:16:19"
$wsZh.Range("L3").Value = ""
$wsZh.Range("M3").Value = "True"
$wsZh.Range("N3").Value = ""
$wsZh.Range("O3").Value = "False"
$wsZh.Range("P3").Value = ""

$wsZh.Hyperlinks.Add($wsZh.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/50f8611c0f7b5f3ff72b035e33299bbcfd184205/e2e/$renamedFile", "", "", $renamedFile) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/1d9cb9ea94d678c4c559803b678a780e457f5a54/e2e/$renamedFile", "", "", $renamedFile) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/50f8611c0f7b5f3ff72b035e33299bbcfd184205/e2e/$newFile", "", "", $newFile) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/1d9cb9ea94d678c4c559803b678a780e457f5a54/e2e/$newFile", "", "", $newFile) | Out-Null

$loZh = $wsZh.ListObjects.Item(1)
$loZh.Resize($wsZh.Range("A1:P3"))

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("A1").Hyperlinks.Delete()

$deRenamedXlf = "$renamedGuid.$renamedHash.de-de.xlf"
$deNewXlf = "$newGuid.$newHash.de-de.xlf"

# Row 2: rename in place, refresh handoff/handback timestamps.
$wsDe.Range("A2").Value = $renamedFile
$wsDe.Range("B2").Value = ".md"
$wsDe.Range("C2").Value = "Handed back: in sync with en-US"
$wsDe.Range("D2").Value = "e2e"
$wsDe.Range("E2").Value = "ht"
$wsDe.Range("F2").Value = "False"
$wsDe.Range("G2").Value = $deRenamedXlf
$wsDe.Range("H2").Value = "2016-08-29 21:15:56"
$wsDe.Range("I2").Value = $renamedFile
$wsDe.Range("J2").Value = $deRenamedXlf
$wsDe.Range("K2").Value = "2016-08-29 21:16:27"
$wsDe.Range("L2").Value = ""
$wsDe.Range("M2").Value = "True"
$wsDe.Range("N2").Value = ""
$wsDe.Range("O2").Value = "False"
$wsDe.Range("P2").Value = ""

# Row 3: brand-new file.
$wsDe.Range("A3").Value = $newFile
$wsDe.Range("B3").Value = ".md"
$wsDe.Range("C3").Value = "Handed back: in sync with en-US"
$wsDe.Range("D3").Value = "e2e"
$wsDe.Range("E3").Value = "ht"
$wsDe.Range("F3").Value = "True"
$wsDe.Range("G3").Value = $deNewXlf
$wsDe.Range("H3").Value = "2016-08-29 21:15:56"
$wsDe.Range("I3").Value = $newFile
$wsDe.Range("J3").Value = $deNewXlf
$wsDe.Range("K3").Value = "2016-08-29 21:16:27"
$wsDe.Range("L3").Value = ""
$wsDe.Range("M3").Value = "True"
$wsDe.Range("N3").Value = ""
$wsDe.Range("O3").Value = "False"
$wsDe.Range("P3").Value = ""

$wsDe.Hyperlinks.Add($wsDe.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/50f8611c0f7b5f3ff72b035e33299bbcfd184205/e2e/$renamedFile", "", "", $renamedFile) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/6fab59b094eb478680e893267b57ae73dcfc0b90/e2e/$renamedFile", "", "", $renamedFile) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/50f8611c0f7b5f3ff72b035e33299bbcfd184205/e2e/$newFile", "", "", $newFile) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/6fab59b094eb478680e893267b57ae73dcfc0b90/e2e/$newFile", "", "", $newFile) | Out-Null

$loDe = $wsDe.ListObjects.Item(1)
$loDe.Resize($wsDe.Range("A1:P3"))

Write-Host "Handback status report refreshed."
